$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add column H ("Ωρες") with the same style as the existing header row ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = 'Ωρες'

# Row 2
$ws.Range("A2").Value = 'ATHENS DENTAL PROJECT - ΜΠΟΥΛΝΤΗΣ ΓΙΩΡΓΟΣ'
$ws.Range("B2").Value = 'Λεωφόρος Βουλιαγμένης 604 & Υμηττού 1, Ελληνικό, 16777, ΑΤΤΙΚΗΣ'
$ws.Range("C2").Value = 'Οδοντιατρικό Κέντρο - Χειρουργός Οδοντίατρος - Γενική Οδοντιατρική - Ορθοδοντική'
$ws.Range("D2").Value = '''2111196707'
$ws.Range("E2").Value = '''6978507450'
$ws.Range("F2").Value = 'http://athensdentalproject.gr'
$ws.Range("G2").Value = 'mailTo:info@athensdentalproject.gr'
$ws.Range("H2").Value = ""

# Row 3
$ws.Range("A3").Value = 'MODERN DENTAL CENTER -ΑΧΙΛΛΑΔΕΛΗΣ ΑΓΓΕΛΟΣ'
$ws.Range("B3").Value = 'Λευκάδος 3, Γλυκά Νερά, 15354, ΑΤΤΙΚΗΣ'
$ws.Range("C3").Value = 'Χειρουργός Οδοντίατρος - Σύγχρονο Οδοντιατρικό Κέντρο – Περιοδοντολογία  Ενδοδοντία - Προσθετική & Επανορθωτική Οδοντιατρική – Εμφυτεύματα - Γναθοπροσωπική Χειρουργική - Ορθοδοντική Παιδιών & Ενηλίκων -Παιδοδοντιατρική'
$ws.Range("D3").Value = '''2106659317'
$ws.Range("E3").Value = '''6974433659'
$ws.Range("F3").Value = 'http://www.moderndental.gr'
$ws.Range("G3").Value = 'mailTo:achilladelis@gmail.com'
$ws.Range("H3").Value = ""

# Row 4
$ws.Range("A4").Value = 'ΠΡΑΣΣΑ ΑΛΕΞΑΝΔΡΟΥ ΑΝΑΣΤΑΣΙΑ - ΜΑΡΙΑ DR. MED'
$ws.Range("B4").Value = 'Λ. Σουνίου 38, ισόγειο, Μαρκόπουλο, 19003, ΑΤΤΙΚΗΣ'
$ws.Range("C4").Value = 'Χειρουργός Οδοντίατρος Εξειδικευθείς Στην Στοματολογία & Περιοδοντολογία'
$ws.Range("D4").Value = '''6936730545'
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 'http://www.stomatomed.gr'
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""

# Row 5
$ws.Range("A5").Value = 'DC ORTHODONTIST - ΧΑΡΑΛΑΜΠΟΠΟΥΛΟΥ ΔΗΜΗΤΡΑ'
$ws.Range("B5").Value = 'Πλαταιών 2, Πλησίον σταθμού ΗΣΑΠ, Μαρούσι, 15124, ΑΤΤΙΚΗΣ'
$ws.Range("C5").Value = 'Ειδικός Ορθοδοντικός'
$ws.Range("D5").Value = '''2108055993'
$ws.Range("E5").Value = '''6932756465'
$ws.Range("F5").Value = 'https://www.dcorthoclinic.com'
$ws.Range("G5").Value = 'mailTo:dcortho@outlook.com'
$ws.Range("H5").Value = ""

# Row 6
$ws.Range("A6").Value = 'ΚΑΚΛΑΜΑΝΗΣ ΑΓΓΕΛΟΣ'
$ws.Range("B6").Value = 'Θηβών 210, Περιστέρι, 12134, ΑΤΤΙΚΗΣ'
$ws.Range("C6").Value = 'Οδοντιατρική Κλινική - Προσθετική - Ενδοδοντία – Γναθοχειρουργική - Ορθοδοντική'
$ws.Range("D6").Value = '''2105730784'
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 'mailTo:kaklamanisdentalcare@gmail.com'
$ws.Range("H6").Value = ""

# Row 7
$ws.Range("A7").Value = 'SPECIALIZED DENTISTRY OF ATHENS - ΛΩΛΗΣ Θ. ΧΡΗΣΤΟΣ'
$ws.Range("B7").Value = 'Ψυχάρη 1 & Στρατήγη, Φάρος, Νέο Ψυχικό, 15451, ΑΤΤΙΚΗΣ'
$ws.Range("C7").Value = 'Αισθητική Οδοντιατρική - Εμφυτευματολογία - Προσθετική - Παιδοδοντία - Περιοδοντολογία - Ορθοδοντική'
$ws.Range("D7").Value = '''2106741600'
$ws.Range("E7").Value = '''6974631121'
$ws.Range("F7").Value = 'http://www.athensdentistry.gr'
$ws.Range("G7").Value = 'mailTo:sda@athensdentistry.gr'
$ws.Range("H7").Value = ""

# Row 8
$ws.Range("A8").Value = 'ΟΔΟΝΤΙΑΤΡΙΚΟ ΚΕΝΤΡΟ ΠΑΛΑΙΟΥ ΦΑΛΗΡΟΥ - ΜΑΡΙΝΑΚΗΣ ΜΑΥΡΟΜΜΑΤΗ ΣΙΑΜΕΤΗ Ο.Ε'
$ws.Range("B8").Value = 'Βενιζέλου Ελευθερίου 186, Παλαιό Φάληρο, 17563, ΑΤΤΙΚΗΣ'
$ws.Range("C8").Value = 'Το Οδοντιατρικό Κέντρο Παρέχει Υπηρεσίες Που Αφορούν Σε Αισθητική Οδοντιατρική - Λεύκανση – Εμφυτεύματα – Προσθετική - Ενδοδοντικές & Περιοδοντικές Θεραπείες - Παιδοδοντία - Ορθοδοντική'
$ws.Range("D8").Value = '''2109819100'
$ws.Range("E8").Value = '''6937303142'
$ws.Range("F8").Value = 'https://dentalclinicpaliofaliro.gr/'
$ws.Range("G8").Value = 'mailTo:dentalclinicps@yahoo.com'
$ws.Range("H8").Value = ""

# Row 9
$ws.Range("A9").Value = 'ΚΟΤΙΝΑΣ ΑΝΑΣΤΑΣΙΟΣ'
$ws.Range("B9").Value = 'Μεσογείων 3, 2ος όροφος, Αθήνα - Αμπελόκηποι, 11526, ΑΤΤΙΚΗΣ'
$ws.Range("C9").Value = 'Ειδικός Ορθοδοντικός'
$ws.Range("D9").Value = '''2107755320'
$ws.Range("E9").Value = '''6944636060'
$ws.Range("F9").Value = 'http://www.akotinas.gr'
$ws.Range("G9").Value = 'mailTo:akotinas@otenet.gr'
$ws.Range("H9").Value = ""

# Row 10
$ws.Range("A10").Value = 'ODOUS MEDICA ΒΑΡΔΑΚΑΣΤΑΝΗ ΑΘΗΝΑ'
$ws.Range("B10").Value = 'Καραϊσκάκη 124, Γλυκά Νερά, , ΑΤΤΙΚΗΣ'
$ws.Range("C10").Value = 'Λεύκανση - Εμφυτεύματα – Ενδοδοντία – Περιοδοντολογία – Προσθετική - Εμφράξεις'
$ws.Range("D10").Value = '''2106041190'
$ws.Range("E10").Value = '''6950568710'
$ws.Range("F10").Value = 'http://www.avardakastani.gr'
$ws.Range("G10").Value = 'mailTo:v.k.athina@hotmail.com'
$ws.Range("H10").Value = ""

# Row 11
$ws.Range("A11").Value = 'ΜΑΚΡΗΣ ΛΕΩΝΙΔΑΣ'
$ws.Range("B11").Value = 'Πανόρμου 119, 2ος όροφος, Αθήνα - Αμπελόκηποι, 11524, ΑΤΤΙΚΗΣ'
$ws.Range("C11").Value = 'Ενδοδοντολόγος - Περιοδοντολόγος - Χειρουργική Στόματος - Χειρουργική Εξαγωγή Φρονιμιτών - Δεκάλεπτη Λεύκανση Δοντιών Με Laser - Καθαρισμός Δοντιών'
$ws.Range("D11").Value = '''6944594547'
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 'http://www.odontiatros-makris.gr'
$ws.Range("G11").Value = 'mailTo:leonmakris@hotmail.com'
$ws.Range("H11").Value = ""

# --- Strip the quote-prefix formatting artifact from forced-text numeric cells, ---
# --- restoring the default (unstyled) cell format while keeping their Text type. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
